$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 2
$ws.Range("E2").Value = 29

# Row 4
$ws.Range("E4").Value = 7

# Row 15
$ws.Range("E15").Value = 98
$ws.Range("F15").Value = 43
$ws.Range("H15").Value = 43

# Row 18
$ws.Range("E18").Value = 52

# Row 25
$ws.Range("E25").Value = 11
$ws.Range("F25").Value = 3
$ws.Range("H25").Value = 3

# Row 26
$ws.Range("E26").Value = 14

# Row 28
$ws.Range("E28").Value = 6

# Row 29
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 2

# Row 42
$ws.Range("E42").Value = 18

# Row 47
$ws.Range("E47").Value = 39
$ws.Range("F47").Value = 21
$ws.Range("H47").Value = 21

# Row 62
$ws.Range("E62").Value = 15

# Row 64
$ws.Range("F64").Value = 12
$ws.Range("H64").Value = 12

# Row 65
$ws.Range("E65").Value = 17
$ws.Range("F65").Value = 6
$ws.Range("H65").Value = 6

# Row 71
$ws.Range("E71").Value = 15

# Row 75
$ws.Range("E75").Value = 9

# Row 79
$ws.Range("E79").Value = 14
